$wb = $excel.ActiveWorkbook

# --- Rename the three sheets, dropping the dot from "hojaejemplo2.N" -> "hojaejemplo2N" ---
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

$ws1.Name = "hojaejemplo21"
$ws2.Name = "hojaejemplo22"
$ws3.Name = "hojaejemplo23"

# --- On the 3rd sheet, mirror column A into column B for rows 1-2 ---
# A1/A2 already hold the shared strings "a"/"b"; writing the same text into
# B1/B2 reuses those same shared-string entries.
$ws3.Range("B1").Value = $ws3.Range("A1").Value()
$ws3.Range("B2").Value = $ws3.Range("A2").Value()

# --- Selections ---
# Sheet 1 keeps its active cell at A5 (selection otherwise unchanged).
$ws1.Select()
$ws1.Range("A5").Select()

# Sheet 2 keeps its active cell at A6 (selection otherwise unchanged).
$ws2.Select()
$ws2.Range("A6").Select()

# Sheet 3 becomes the active sheet/tab, with B1:B2 selected (active cell B1).
$ws3.Select()
$ws3.Range("B1:B2").Select()
